# Chaitanya basic scenarios add
# Updates the "Summary", "Repayment Schedule" and "Transactions" sheets:
#  - Summary!F3 value changes, and the selected cell moves to B4
#  - Repayment Schedule selection moves to C5
#  - Transactions!A2 / A3 values change, and the selection moves to C2
#    (Transactions is left as the active/selected tab, matching the source)

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("Summary")
$summary.Range("F3").Value = 1026.54
$summary.Range("B4").Select()

$repaymentSchedule = $wb.Worksheets.Item("Repayment Schedule")
$repaymentSchedule.Range("C5").Select()

$transactions = $wb.Worksheets.Item("Transactions")
$transactions.Range("A2").Value = 194
$transactions.Range("A3").Value = 193
$transactions.Range("C2").Select()
